$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.925.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.359.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.25%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -11.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.358.46"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.05%  "
$ws.Range("E10").Value = "  -2.39%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("E12").Value = "  -3.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.340"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.43%  "
$ws.Range("E14").Value = "  -4.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.779.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.766.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.352.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.63%  "
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("E22").Value = "  -6.86%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.471.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0892"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "504.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.89%  "
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("E34").Value = "  -5.89%  "
$ws.Range("E35").Value = "  -3.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.78%  "
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.35%  "
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "138.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.58%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "138.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.08%  "
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("E48").Value = "  -4.80%  "
$ws.Range("E49").Value = "  -9.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.571"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.04%  "
$ws.Range("E51").Value = "  -2.20%  "
